# Update the "Wheels" CANBus node-status row (row 15) so each short code is
# followed by its expanded/clarified description, per the commit:
#   "Update TPS and BMS code for Adruino, update BLE Services.xlsx"
#
# Old layout (F15:S15): bms, tps, sas, whl, imu, int, flw, frw, rlw, rrw, fll, frl, rll, rrl
# New layout (F15:O15): bms (...), tps (...), sas (...), imu (...), fw (...),
#                        rlw (...), rrw (...), fl (...), rl (...), int (...)
# Columns P15:S15 are no longer used and are cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F15").Value = "bms (battery management system)"
$ws.Range("G15").Value = "tps (throttle and brake position sensor)"
$ws.Range("H15").Value = "sas (steering angle sensor)"
$ws.Range("I15").Value = "imu (inertial measurment unit)"
$ws.Range("J15").Value = "fw (front wheels)"
$ws.Range("K15").Value = "rlw (rear left wheel)"
$ws.Range("L15").Value = "rrw (rear right wheel)"
$ws.Range("M15").Value = "fl (front light)"
$ws.Range("N15").Value = "rl (rear light)"
$ws.Range("O15").Value = "int (interior light)"

# I15 no longer wraps text (moves to plain/default style); K15 now does.
$ws.Range("I15").Style = "Normal"
$ws.Range("K15").WrapText = $true

# The trailing four cells that used to hold fll/frl/rll/rrl are now blank.
$ws.Range("P15:S15").ClearContents()

# Reflect the author's final view/cursor position when they saved the file:
# scrolled right so column E is the leftmost visible column, with F16 selected.
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1
$ws.Range("F16").Select()
